$d = $word.ActiveDocument

$replacements = @(
    @("291×3=", "992×4="),
    @("730×8=", "383×4="),
    @("931×8=", "174×3="),
    @("137×7=", "169×4="),
    @("982×3=", "887×8="),
    @("493×5=", "459×2="),
    @("452×5=", "677×8="),
    @("665×8=", "272×9="),
    @("289×2=", "789×7="),
    @("461×2=", "943×2="),
    @("472×8=", "835×9="),
    @("995×6=", "264×5="),
    @("894×5=", "401×4="),
    @("866×3=", "809×6="),
    @("299×2=", "536×5="),
    @("178×8=", "933×9="),
    @("945×3=", "473×5="),
    @("421×7=", "136×9="),
    @("106×4=", "207×5="),
    @("266×2=", "837×8="),
    @("675×6=", "716×9="),
    @("499×8=", "706×7="),
    @("814×3=", "573×6="),
    @("673×9=", "269×7="),
    @("490×3=", "942×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
